$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New tender rows (2-9) replacing the old ones
$data = @(
    @("Услуги по перезарядке и ремонту огнетушителей, согласно технического задания (приложение 1 к государственному контракту)", "01:27:03", "ФКУ ИК-6 УФСИН РОССИИ ПО Г. САНКТ-ПЕТЕРБУРГУ И ЛЕНИНГРАДСКОЙ ОБЛАСТИ", "66 183,43 руб.", "https://agregatoreat.ru/purchase/1392493/order-info"),
    @("Заправка огнетушителей, 4 шт", "01:14:23", "ОТДЕЛ МВД РОССИИ ПО ИГЛИНСКОМУ РАЙОНУ", "840,00 руб.", "https://agregatoreat.ru/purchase/1392389/order-info"),
    @("Переосвидетельствование и перезарядка огнетушителей., 30 шт", "23:01:43", "ОМВД РОССИИ ПО КАМЫЗЯКСКОМУ РАЙОНУ", "3 600,00 руб.", "https://agregatoreat.ru/purchase/1392302/order-info"),
    @("Перезарядка огнетушителя ОП-2, 4 шт", "22:18:04", "ОТДЕЛ ГФС РОССИИ В Г. ВЕЛИКОМ НОВГОРОДЕ", "1 900,00 руб.", "https://agregatoreat.ru/purchase/1392001/order-info"),
    @("Одежда и обувь, Автомобильные товары и запчасти", "20:42:01", "ФКУКИИ ""КУЛЬТУРНЫЙ ЦЕНТР МВД РОССИИ""", "7 500,00 руб.", "https://agregatoreat.ru/purchase/1391454/order-info"),
    @("Оказание услуг по вывозу и утилизации огнетушителей, 1 усл. ед", "17:19:56", "УТ МВД РОССИИ ПО ДФО", "29 200,00 руб.", "https://agregatoreat.ru/purchase/1389772/order-info"),
    @("Выполнение работ по техническому обслуживанию углекислотных огнетушителей, 300 шт", "01:39:34", "ФСО РОССИИ", "195 000,00 руб.", "https://agregatoreat.ru/purchase/1388261/order-info"),
    @("Техническое обслуживание и перезарядка огнетушителей, 132 шт", "01:01:09", "ФГКУ КОМБИНАТ ""СИГНАЛ"" ИМ. А.А. ГРИГОРЬЕВА РОСРЕЗЕРВА", "95 850,00 руб.", "https://agregatoreat.ru/purchase/1388110/order-info")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}

# Remove the now-unused last row (old row 10)
$ws.Rows.Item(10).Delete()
